# "Added crdc login backup code"
# The previously-active backup codes (rows 3, 4 and 11) have been consumed,
# so they are removed from the sheet. The next unused code moves into the
# "current" slot (A2), and the last remaining unused code (A12) stays put.
# Selection moves to A2, the new current backup code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 becomes the new current backup code
$ws.Range("A2").Value = "NHYK5008HQDA"

# The consumed backup codes are cleared out
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()
$ws.Range("A11").ClearContents()

# A12 keeps its existing (still-unused) backup code "ZADGNDVPP03M"

# Put the selection on the new current backup code
$ws.Range("A2").Select() | Out-Null
